# Fixed table lines and handled file open exception
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (first data row) - update sample/demo values
$ws.Range("M3").Value = "NewUsedEquipment 1"
$ws.Range("R3").Value = "1. 100 2. 200"
$ws.Range("AO3").Value = "1. 22 2. 33"
$ws.Range("AP3").Value = "1. Description^p 2. ofGoodss"
$ws.Range("AQ3").Value = "1. 1 2. 5"

# Row 4 (second data row)
$ws.Range("M4").Value = "NNN"
$ws.Range("AO4").Value = "cthno"

# Row 5 (third data row)
$ws.Range("M5").Value = "NewUsedEquipment 313 !@#$%^&*()_+"
